{"js": "// Fix the \"Plant *\" summary tables:\n//   1. Column header text \"std.err\" -> \"se\" (one per table, 8 tables total).\n//   2. The corresponding (5th) table-grid column width 1071 -> 971 twips,\n//      so it lines up with the other narrow (971-wide) columns.\n\n// --- 1. Rename every \"std.err\" header cell to \"se\" -----------------------\nconst stdErrRanges = context.document.body.search(\"std.err\", { matchCase: true });\nstdErrRanges.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < stdErrRanges.items.length; i++) {\n  stdErrRanges.items[i].insertText(\"se\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2. Narrow the 5th column (index 4) of every table from 1071 -> 971 --\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst TWIPS_PER_POINT = 20;\nconst targetWidthTwips = 971;\nconst STD_ERR_COLUMN_INDEX = 4; // 0-based: term, est, ci.lb, ci.ub, std.err/se\n\nfor (let i = 0; i < tables.items.length; i++) {\n  const table = tables.items[i];\n  const headerCell = table.getCell(0, STD_ERR_COLUMN_INDEX);\n  headerCell.columnWidth = targetWidthTwips / TWIPS_PER_POINT;\n}\nawait context.sync();\n", "ps1": "# Fix the \"Plant *\" summary tables:\n#   1. Column header text \"std.err\" -> \"se\" (one per table, 8 tables total).\n#   2. The corresponding (5th) table-grid column width 1071 -> 971 twips,\n#      so it lines up with the other narrow (971-wide) columns.\n\n$d = $word.ActiveDocument\n\n# --- 1. Rename every \"std.err\" header cell to \"se\" ------------------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"std.err\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 1\nwhile ($rng.Find.Execute()) {\n    $rng.Text = \"se\"\n}\n\n# --- 2. Narrow the 5th column (index 4 / cell 5) of every table -----------\n$targetWidthTwips = 971\n$targetWidthPoints = $targetWidthTwips / 20\n$stdErrColumn = 5   # 1-based: term, est, ci.lb, ci.ub, std.err/se\n\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    $cell = $t.Cell(1, $stdErrColumn)\n    $cell.Width = $targetWidthPoints\n}\n"}
